$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove all existing values/formatting in the used range
$ws.UsedRange.Clear()

# Row 1: new "Decision" / "Hello World" prompt header
$ws.Range("A1").Value = "Decision"
$ws.Range("B1").Value = "Hello World"
$ws.Range("C1").ShrinkToFit = $false
$ws.Range("D1").ShrinkToFit = $false

# Row 2: Options header + the two entries (Banana Firm / Pear Company)
$ws.Range("A2").Value = "Options"
$ws.Range("B2").Value = "Banana Firm"
$ws.Range("C2").Value = "Pear Company"
$ws.Range("D2").ShrinkToFit = $false

# Row 3: thin spacer row
$ws.Rows(3).RowHeight = 4.5
$ws.Range("D3").ShrinkToFit = $false

# Row 4: Criteria header + first criterion (Entry / Difficulty)
$ws.Range("A4").Value = "Criteria"
$ws.Range("A4").ShrinkToFit = $false
$ws.Range("B4").Value = "Entry"
$ws.Range("C4").Value = "Difficulty"

# Row 5: second criterion (Time)
$ws.Range("C5").Value = "Time"

# Row 6: third criterion (Content / Work hour / Length)
$ws.Range("B6").Value = "Content"
$ws.Range("C6").Value = "Work hour"
$ws.Range("D6").Value = "Length"

# Row 7: fourth criterion (Sleep schedule)
$ws.Range("D7").Value = "Sleep schedule"

# Column widths: B, C, D now share the wider 16.53125-ish width
$ws.Columns("B").ColumnWidth = 15.7
$ws.Columns("C").ColumnWidth = 15.7
$ws.Columns("D").ColumnWidth = 15.7

# Update selection to match the saved workbook state
$ws.Range("C9").Select()
